# Apply crypto price/volume updates for Sun Mar 3 2024 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D3").NumberFormat = "@"
$ws.Range("D5:D14").NumberFormat = "@"
$ws.Range("D16:D18").NumberFormat = "@"
$ws.Range("D20:D24").NumberFormat = "@"
$ws.Range("D26:D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38:D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "62.256.18"
$ws.Range("E2").Value = "  +0.67%  "

# Row 3
$ws.Range("D3").Value = "3.427.79"
$ws.Range("E3").Value = "  +0.58%  "

# Row 4
$ws.Range("E4").Value = "  +0.43%  "

# Row 5
$ws.Range("D5").Value = "413.13"
$ws.Range("E5").Value = "  +1.14%  "

# Row 6
$ws.Range("D6").Value = "129.30"
$ws.Range("E6").Value = "  +0.85%  "

# Row 7
$ws.Range("D7").Value = "0.621"
$ws.Range("E7").Value = "  -2.02%  "

# Row 8
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.06%  "

# Row 9
$ws.Range("D9").Value = "0.725"
$ws.Range("E9").Value = "  -0.60%  "

# Row 10
$ws.Range("D10").Value = "0.139"
$ws.Range("E10").Value = "  +0.97%  "

# Row 11
$ws.Range("D11").Value = "42.85"
$ws.Range("E11").Value = "  +1.02%  "

# Row 12
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "9.19"
$ws.Range("E12").Value = "  +1.35%  "

# Row 13
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "3.979.76"
$ws.Range("E13").Value = "  +0.98%  "

# Row 14
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").Value = "0.0000214"
$ws.Range("E14").Value = "  +7.09%  "

# Row 15
$ws.Range("E15").Value = "  -0.39%  "

# Row 16
$ws.Range("D16").Value = "20.43"
$ws.Range("E16").Value = "  -3.60%  "

# Row 17
$ws.Range("D17").Value = "3.394.61"
$ws.Range("E17").Value = "  -0.33%  "

# Row 18
$ws.Range("D18").Value = "12.59"
$ws.Range("E18").Value = "  +4.57%  "

# Row 19
$ws.Range("E19").Value = "  -0.30%  "

# Row 20
$ws.Range("D20").Value = "62.327.10"

# Row 21
$ws.Range("D21").Value = "470.46"
$ws.Range("E21").Value = "  +4.66%  "

# Row 22
$ws.Range("D22").Value = "91.15"
$ws.Range("E22").Value = "  -0.61%  "

# Row 23
$ws.Range("D23").Value = "3.25"
$ws.Range("E23").Value = "  +2.88%  "

# Row 24
$ws.Range("D24").Value = "13.22"
$ws.Range("E24").Value = "  +2.58%  "

# Row 25
$ws.Range("E25").Value = "  +2.02%  "

# Row 26
$ws.Range("D26").Value = "9.93"
$ws.Range("E26").Value = "  +14.13%  "

# Row 27
$ws.Range("D27").Value = "33.02"
$ws.Range("E27").Value = "  -1.12%  "

# Row 28
$ws.Range("E28").Value = "  +0.63%  "

# Row 29
$ws.Range("D29").Value = "7.74"
$ws.Range("E29").Value = "  +2.54%  "

# Row 30
$ws.Range("E30").Value = "  -0.84%  "

# Row 31
$ws.Range("D31").Value = "2.62"
$ws.Range("E31").Value = "  -4.06%  "

# Row 32
$ws.Range("E32").Value = "  -1.04%  "

# Row 33
$ws.Range("E33").Value = "  -1.95%  "

# Row 34
$ws.Range("D34").Value = "41.07"
$ws.Range("E34").Value = "  -4.12%  "

# Row 35
$ws.Range("E35").Value = "  +0.02%  "

# Row 36
$ws.Range("D36").Value = "57.61"
$ws.Range("E36").Value = "  +8.33%  "

# Row 38
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.30%  "

# Row 39
$ws.Range("D39").Value = "3.05"
$ws.Range("E39").Value = "  +5.09%  "

# Row 40
$ws.Range("D40").Value = "0.327"
$ws.Range("E40").Value = "  +4.31%  "

# Row 41
$ws.Range("D41").Value = "0.134"
$ws.Range("E41").Value = "  +0.26%  "

# Row 42
$ws.Range("E42").Value = "  -0.91%  "

# Row 43
$ws.Range("D43").Value = "144.47"
$ws.Range("E43").Value = "  +2.28%  "

# Row 44
$ws.Range("E44").Value = "  +9.87%  "

# Row 45
$ws.Range("D45").Value = "2.07"
$ws.Range("E45").Value = "  +4.96%  "

# Row 46
$ws.Range("D46").Value = "4.33"
$ws.Range("E46").Value = "  +3.41%  "

# Row 47
$ws.Range("D47").Value = "2.42"
$ws.Range("E47").Value = "  +19.54%  "

# Row 48
$ws.Range("D48").Value = "16.47"
$ws.Range("E48").Value = "  -0.23%  "

# Row 49
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "22.15"
$ws.Range("E49").Value = "  -0.43%  "

# Row 50
$ws.Range("B50").Value = "PEPE"
$ws.Range("C50").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D50").Value = "0.0₃0517"
$ws.Range("E50").Value = "  +27.71%  "

# Row 51
$ws.Range("D51").Value = "111.14"
$ws.Range("E51").Value = "  +4.91%  "
